$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Remove the duplicate "Contact" row (row 11 duplicated row 10's content).
$ws.Rows.Item(11).Delete()

# Version bump 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date refresh
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Row 9 was "Publisher" with an empty value -> now carries the publisher name.
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 used to be the (now-removed) duplicate "Contact" row -> becomes "Jurisdiction".
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# "Case Sensitive" row now carries a "true" text value. Building it from a
# formula and then freezing it to a value keeps it a genuine text cell
# (a literal "true" would be auto-coerced to the boolean TRUE by Excel).
$caseSensitiveCell = $ws.Range("B14")
$caseSensitiveCell.Formula = '="tru"&"e"'
$caseSensitiveCell.Copy() | Out-Null
$caseSensitiveCell.PasteSpecial(-4163) | Out-Null  # xlPasteValues
$excel.CutCopyMode = $false
